$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Apply updated crypto price/volume data (and the Huobi/MX row swap)

$ws.Range("D2").Value = "'26.165.36"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "'  -4.31%  "
$ws.Range("E2").Style = "Normal"

$ws.Range("D3").Value = "'1.654.68"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "'  -3.50%  "
$ws.Range("E3").Style = "Normal"

$ws.Range("D4").Value = "'1.011"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "'  +0.29%  "
$ws.Range("E4").Style = "Normal"

$ws.Range("D5").Value = "'216.36"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "'  -3.74%  "
$ws.Range("E5").Style = "Normal"

$ws.Range("D6").Value = "'0.5123"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "'  -2.94%  "
$ws.Range("E6").Style = "Normal"

$ws.Range("E7").Value = "'  +0.45%  "
$ws.Range("E7").Style = "Normal"

$ws.Range("D8").Value = "'0.2597"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "'  -1.79%  "
$ws.Range("E8").Style = "Normal"

$ws.Range("D9").Value = "'0.06443"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "'  -3.34%  "
$ws.Range("E9").Style = "Normal"

$ws.Range("D10").Value = "'19.76"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "'  -4.85%  "
$ws.Range("E10").Style = "Normal"

$ws.Range("D11").Value = "'0.07801"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "'  +0.57%  "
$ws.Range("E11").Style = "Normal"

$ws.Range("D12").Value = "'1.648.69"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "'  -3.71%  "
$ws.Range("E12").Style = "Normal"

$ws.Range("D13").Value = "'4.283"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "'  -4.17%  "
$ws.Range("E13").Style = "Normal"

$ws.Range("D14").Value = "'1.883.18"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "'  -3.45%  "
$ws.Range("E14").Style = "Normal"

$ws.Range("D15").Value = "'0.5488"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "'  -5.28%  "
$ws.Range("E15").Style = "Normal"

$ws.Range("D16").Value = "'0.0₅8007"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "'  -2.08%  "
$ws.Range("E16").Style = "Normal"

$ws.Range("D17").Value = "'63.98"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "'  -5.53%  "
$ws.Range("E17").Style = "Normal"

$ws.Range("D18").Value = "'26.194.30"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "'  -4.25%  "
$ws.Range("E18").Style = "Normal"

$ws.Range("D19").Value = "'1.011"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "'  +0.25%  "
$ws.Range("E19").Style = "Normal"

$ws.Range("D20").Value = "'208.24"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "'  -5.17%  "
$ws.Range("E20").Style = "Normal"

$ws.Range("D21").Value = "'4.400"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "'  -5.35%  "
$ws.Range("E21").Style = "Normal"

$ws.Range("D22").Value = "'10.09"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "'  -3.16%  "
$ws.Range("E22").Style = "Normal"

$ws.Range("D23").Value = "'6.039"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "'  +0.09%  "
$ws.Range("E23").Style = "Normal"

$ws.Range("D24").Value = "'1.012"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "'  +0.42%  "
$ws.Range("E24").Style = "Normal"

$ws.Range("E25").Value = "'  +7.99%  "
$ws.Range("E25").Style = "Normal"

$ws.Range("D26").Value = "'144.36"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "'  -0.55%  "
$ws.Range("E26").Style = "Normal"

$ws.Range("D27").Value = "'0.1170"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "'  -2.87%  "
$ws.Range("E27").Style = "Normal"

$ws.Range("D28").Value = "'6.956"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "'  -3.69%  "
$ws.Range("E28").Style = "Normal"

$ws.Range("D29").Value = "'15.82"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "'  -2.24%  "
$ws.Range("E29").Style = "Normal"

$ws.Range("D30").Value = "'0.05085"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "'  -5.04%  "
$ws.Range("E30").Style = "Normal"

$ws.Range("D31").Value = "'1.244"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "'  -3.89%  "
$ws.Range("E31").Style = "Normal"

$ws.Range("D32").Value = "'3.346"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "'  -3.89%  "
$ws.Range("E32").Style = "Normal"

$ws.Range("E33").Value = "'  -4.33%  "
$ws.Range("E33").Style = "Normal"

$ws.Range("D34").Value = "'1.553"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "'  -5.05%  "
$ws.Range("E34").Style = "Normal"

$ws.Range("B35").Value = "'HuobiToken"
$ws.Range("B35").Style = "Normal"
$ws.Range("C35").Value = "'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
$ws.Range("C35").Style = "Normal"
$ws.Range("D35").Value = "'2.362"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "'  -1.62%  "
$ws.Range("E35").Style = "Normal"

$ws.Range("B36").Value = "'MXToken"
$ws.Range("B36").Style = "Normal"
$ws.Range("C36").Value = "'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("C36").Style = "Normal"
$ws.Range("D36").Value = "'2.724"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "'  -4.36%  "
$ws.Range("E36").Style = "Normal"

$ws.Range("D37").Value = "'0.9190"
$ws.Range("D37").Style = "Normal"

$ws.Range("D38").Value = "'1.174.03"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "'  +1.26%  "
$ws.Range("E38").Style = "Normal"

$ws.Range("D39").Value = "'0.5705"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "'  -2.93%  "
$ws.Range("E39").Style = "Normal"

$ws.Range("D40").Value = "'0.01584"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "'  -3.98%  "
$ws.Range("E40").Style = "Normal"

$ws.Range("E41").Value = "'  +0.41%  "
$ws.Range("E41").Style = "Normal"

$ws.Range("D42").Value = "'2.567"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "'  -0.44%  "
$ws.Range("E42").Style = "Normal"

$ws.Range("D43").Value = "'5.665"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "'  -2.51%  "
$ws.Range("E43").Style = "Normal"

$ws.Range("D44").Value = "'0.8269"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "'  -1.44%  "
$ws.Range("E44").Style = "Normal"

$ws.Range("D45").Value = "'100.48"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "'  -0.60%  "
$ws.Range("E45").Style = "Normal"

$ws.Range("D46").Value = "'1.795.15"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "'  -3.35%  "
$ws.Range("E46").Style = "Normal"

$ws.Range("D47").Value = "'0.0₈113"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "'  -4.11%  "
$ws.Range("E47").Style = "Normal"

$ws.Range("D48").Value = "'0.4556"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "'  +0.17%  "
$ws.Range("E48").Style = "Normal"

$ws.Range("D49").Value = "'1.009"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "'  +0.02%  "
$ws.Range("E49").Style = "Normal"

$ws.Range("D50").Value = "'55.34"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "'  -3.59%  "
$ws.Range("E50").Style = "Normal"

$ws.Range("D51").Value = "'7.847"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "'  -3.47%  "
$ws.Range("E51").Style = "Normal"
